$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) "Projekthandbuch" -> "Beginn Projekthandbuch" on the "Richtungspfeil 5" shape
#    (second paragraph, single run)
$arrow = $s.Shapes.Item(3)
$arrow.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "Beginn Projekthandbuch"

# 2) "Fertigstellung des Projekts" textbox: move/resize + bump first run font size
$box = $s.Shapes.Item(11)
$box.Top = 280.845511811
$box.Height = 59.3742519685
$box.TextFrame.TextRange.Paragraphs(1).Runs(1).Font.Size = 16
